$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (pushes the existing data rows down by one,
# so old rows 2-6 become rows 3-7)
$ws.Rows.Item(2).Insert()

# Copy formatting (the thin-border cell style) from the row that is now
# row 3 onto the newly inserted (blank) row 2, so it matches the rest of
# the table
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# Plain text values can be assigned directly without Excel re-typing them
$ws.Range("B2").Value = "BEBELAC1+ MDU KLG800"
$ws.Range("C2").Value = "ECPRM1"
$ws.Range("F2").Value = "RT,(E-2B)"

# Numeric-looking values must stay text (matching the rest of the sheet),
# so stage each one in a scratch cell as a text formula result, then
# paste only the value across - this avoids Excel's automatic
# text->number conversion and keeps the style table untouched
$scratch = $ws.Range("Z100")

$scratch.Formula = '="20035632"'
$scratch.Copy()
$ws.Range("A2").PasteSpecial(-4163)

$scratch.Formula = '="1"'
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)

$scratch.Formula = '="1"'
$scratch.Copy()
$ws.Range("E2").PasteSpecial(-4163)

$scratch.Clear()
